$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2:B118 epoch accuracy values (re-run of the training notebook)
$values = @(
    0.734375,
    0.59375,
    0.515625,
    0.5,
    0.375,
    0.359375,
    0.296875,
    0.328125,
    0.28125,
    0.234375,
    0.296875,
    0.375,
    0.34375,
    0.28125,
    0.3125,
    0.203125,
    0.234375,
    0.21875,
    0.25,
    0.25,
    0.25,
    0.25,
    0.25,
    0.25,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.21875,
    0.21875,
    0.21875,
    0.21875,
    0.21875,
    0.21875,
    0.21875,
    0.21875,
    0.21875,
    0.21875,
    0.21875,
    0.21875,
    0.21875,
    0.21875,
    0.21875,
    0.21875,
    0.21875,
    0.21875,
    0.21875,
    0.21875,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.234375,
    0.171875,
    0.265625,
    0.21875,
    0.15625,
    0.15625,
    0.25,
    0.15625,
    0.234375,
    0.265625,
    0.203125,
    0.171875,
    0.265625,
    0.25,
    0.265625,
    0.1639344262295082
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Update the DisplayOutputs repr (new object memory address) in column A rows 102-118
$newAddr = "<__main__.DisplayOutputs object at 0x7fe8c08134f0>"
for ($row = 102; $row -le 118; $row++) {
    $ws.Cells.Item($row, 1).Value = $newAddr
}

# Move the active cell / selection recorded in the sheet view to O19
$ws.Range("A1:XFD1048576").Select()
$ws.Range("O19").Activate()
